$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match-by-match stats (runs, balls, fours, sixes) got reshuffled across rows.
# Row 10 is untouched; every other row (2-9, 11-16) takes on the C:F values
# that used to belong to a different row. Values are text-typed numbers, so
# force the number format to Text ("@") before assigning, which keeps the
# cells stored the same way (text) rather than letting Excel auto-detect them
# as numbers.
$newValues = @{
    2  = @("21", "14", "3", "0")
    3  = @("51", "38", "3", "2")
    4  = @("6",  "7",  "0", "0")
    5  = @("1",  "2",  "0", "0")
    6  = @("3",  "5",  "0", "0")
    7  = @("29", "21", "5", "0")
    8  = @("83", "47", "4", "8")
    9  = @("4",  "3",  "1", "0")
    11 = @("24", "21", "3", "1")
    12 = @("26", "19", "3", "1")
    13 = @("44", "31", "4", "1")
    14 = @("15", "29", "0", "0")
    15 = @("34", "33", "3", "1")
    16 = @("30", "19", "4", "1")
}

$cols = @("C", "D", "E", "F")

foreach ($row in $newValues.Keys) {
    $vals = $newValues[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cellRange = $ws.Range("$($cols[$i])$row")
        $cellRange.NumberFormat = "@"
        $cellRange.Value = $vals[$i]
    }
}
